# Monte Carlo LCA calculation script initial
# Applies the "with multiprocessing" comparison block (columns O:Q) to the
# "For dependant MC" worksheet, plus updated creation-time measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For dependant MC")

# --- Updated "Creation time [sec]" measurements (rows 8 and 12) ---
$ws.Range("D8").Formula = "=1.7375"
$ws.Range("D12").Formula = "=0.84333"
$ws.Range("F12").Value = 0.35

# --- New "number of processes" parameter ---
$ws.Range("P20").Value = 4

# --- New header row 21: "With " <P20> "processes" ---
$ws.Range("O21").Value = "With "
$ws.Range("P21").Formula = '=$P$20'
$ws.Range("Q21").Value = "processes"

# --- New header row 22: same label as L22 ---
$ws.Range("O22").Value = $ws.Range("L22").Value2

# --- New header row 23: same labels as L23:N23 ---
$ws.Range("O23").Value = $ws.Range("L23").Value2
$ws.Range("P23").Value = $ws.Range("M23").Value2
$ws.Range("Q23").Value = $ws.Range("N23").Value2

# --- New per-row time-with-multiprocessing columns (rows 24-27) ---
# (ClearFormats avoids inheriting the referenced cells' number format, so the
#  new cells stay on the default style exactly like the source workbook.)
$ws.Range("O24").Formula = '=(L24)/($P$20)'
$ws.Range("O24").ClearFormats()
$ws.Range("P24").Formula = '=(M24)/($P$20)'
$ws.Range("P24").ClearFormats()
$ws.Range("Q24").Formula = '=(N24)/($P$20)'
$ws.Range("Q24").ClearFormats()

$ws.Range("O25").Formula = '=(L25)/($P$20)'
$ws.Range("O25").ClearFormats()
$ws.Range("P25").Formula = '=(M25)/($P$20)'
$ws.Range("P25").ClearFormats()
$ws.Range("Q25").Formula = '=(N25)/($P$20)'
$ws.Range("Q25").ClearFormats()

$ws.Range("O26").Formula = '=(L26)/($P$20)'
$ws.Range("O26").ClearFormats()
$ws.Range("P26").Formula = '=(M26)/($P$20)'
$ws.Range("P26").ClearFormats()
$ws.Range("Q26").Formula = '=(N26)/($P$20)'
$ws.Range("Q26").ClearFormats()

$ws.Range("O27").Formula = '=(L27)/($P$20)'
$ws.Range("O27").ClearFormats()
$ws.Range("P27").Formula = '=(M27)/($P$20)'
$ws.Range("P27").ClearFormats()
$ws.Range("Q27").Formula = '=(N27)/($P$20)'
$ws.Range("Q27").ClearFormats()

# --- Totals row 28, matching the bold styling used by L28:N28 ---
$ws.Range("O28").Formula = "=SUM(O24:O27)"
$ws.Range("P28").Formula = "=SUM(P24:P27)"
$ws.Range("Q28").Formula = "=SUM(Q24:Q27)"

$ws.Range("O28").ClearFormats()
$ws.Range("P28").ClearFormats()
$ws.Range("Q28").ClearFormats()
$ws.Range("O28").Font.Bold = $true
$ws.Range("P28").Font.Bold = $true
$ws.Range("Q28").Font.Bold = $true
$ws.Range("P28").NumberFormat = "0.00"
$ws.Range("Q28").NumberFormat = "0.00"

# --- Selection matches the author's final cursor position ---
$ws.Range("R26").Select() | Out-Null
